$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8 ("拿走 / nv / -") to make room for
# a new vocabulary entry: 被 / passive / 虛詞 (pushes old rows 8-13 down to 9-14).
$ws.Rows.Item(8).Insert()

# Populate the new row 8.
$ws.Range("A8").Value = "被"
$ws.Range("B8").Value = "passive"
$ws.Range("C8").Value = "虛詞"

# Match the row height used by the other "highlight box" rows (7 & the old/
# new row 8) in this sheet.
$ws.Rows.Item(8).RowHeight = 16.5

# Re-create the red "highlight box" formatting used elsewhere in the sheet
# (see row 8/9 "拿走 nv -") for this newly inserted row: red font + a medium
# red border that reads as one continuous box across A8:C8.
$a8 = $ws.Range("A8")
$a8.Font.Color = 255
$a8.Borders.Item(7).Weight = -4138   # xlEdgeLeft
$a8.Borders.Item(7).Color = 255
$a8.Borders.Item(8).Weight = -4138   # xlEdgeTop
$a8.Borders.Item(8).Color = 255
$a8.Borders.Item(9).Weight = -4138   # xlEdgeBottom
$a8.Borders.Item(9).Color = 255

$b8 = $ws.Range("B8")
$b8.Font.Color = 255
$b8.Borders.Item(8).Weight = -4138   # xlEdgeTop
$b8.Borders.Item(8).Color = 255
$b8.Borders.Item(9).Weight = -4138   # xlEdgeBottom
$b8.Borders.Item(9).Color = 255

$c8 = $ws.Range("C8")
$c8.Font.Color = 255
$c8.Borders.Item(8).Weight = -4138   # xlEdgeTop
$c8.Borders.Item(8).Color = 255
$c8.Borders.Item(9).Weight = -4138   # xlEdgeBottom
$c8.Borders.Item(9).Color = 255
$c8.Borders.Item(10).Weight = -4138  # xlEdgeRight
$c8.Borders.Item(10).Color = 255

# Match the selection left behind by the edit (the newly inserted row).
$ws.Range("A8:C8").Select()

Write-Output "done"
